$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 16, shifting rows 16:48 down to 17:49
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the "flux/depth" case
$ws.Range("A16").Value = "New 200m transport,  flux/depth, dt = 1 hr, j = 2"
$ws.Range("B16").Value = -4.5010000000000003
$ws.Range("C16").Value = -1.5115000000000001
$ws.Range("D16").Value = -1.7454000000000001
$ws.Range("E16").Value = -2.09
$ws.Range("F16").Value = 0.9445
$ws.Range("G16").Value = -0.6676
$ws.Range("H16").Value = 0.4495
$ws.Range("I16").Value = -0.1013

# Update the sheet view: scrolled to the top, selection on H16
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("H16").Select()
